$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2:T11").ClearContents()
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il1rn"
$ws.Range("C2").Value = "Il1r1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("D4").Value = "M1"
$ws.Range("D5").Value = "M2"
$ws.Range("D6").Value = "sCs"
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il1rn"
$ws.Range("C2").Value = "Il1r1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.4115716666666667
$ws.Range("H2").Value = 1.234715
$ws.Range("I2").Value = 0.002110162210096788
$ws.Range("J2").Value = 0.002110162210096788
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 17.79587766666667
$ws.Range("N2").Value = 53.38763299999999
$ws.Range("O2").Value = 0.1793479316144739
$ws.Range("P2").Value = 0.179347931614474
$ws.Range("Q2").Value = 7.324279031066111
$ws.Range("R2").Value = 65.91851127959499
$ws.Range("S2").Value = 0.0003784532277518859
$ws.Range("T2").Value = 0.000378453227751886
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Il1rn"
$ws.Range("C3").Value = "Il1r1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.4115716666666667
$ws.Range("H3").Value = 1.234715
$ws.Range("I3").Value = 0.002110162210096788
$ws.Range("J3").Value = 0.002110162210096788
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 72.39518466666667
$ws.Range("N3").Value = 217.185554
$ws.Range("O3").Value = 0.7296030503252249
$ws.Range("P3").Value = 0.7296030503252251
$ws.Range("Q3").Value = 29.79580681190111
$ws.Range("R3").Value = 268.16226130711
$ws.Range("S3").Value = 0.001539580785167635
$ws.Range("T3").Value = 0.001539580785167635
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Il1rn"
$ws.Range("C4").Value = "Il1r1"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.4115716666666667
$ws.Range("H4").Value = 1.234715
$ws.Range("I4").Value = 0.002110162210096788
$ws.Range("J4").Value = 0.002110162210096788
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.257112
$ws.Range("N4").Value = 0.771336
$ws.Range("O4").Value = 0.00259119028895291
$ws.Range("P4").Value = 0.00259119028895291
$ws.Range("Q4").Value = 0.10582001436
$ws.Range("R4").Value = 0.95238012924
$ws.Range("S4").Value = 0.000005467831826918206
$ws.Range("T4").Value = 0.000005467831826918207
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Il1rn"
$ws.Range("C5").Value = "Il1r1"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4115716666666667
$ws.Range("H5").Value = 1.234715
$ws.Range("I5").Value = 0.002110162210096788
$ws.Range("J5").Value = 0.002110162210096788
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4631083333333333
$ws.Range("N5").Value = 1.389325
$ws.Range("O5").Value = 0.004667233797203165
$ws.Range("P5").Value = 0.004667233797203166
$ws.Range("Q5").Value = 0.1906022685972222
$ws.Range("R5").Value = 1.715420417375
$ws.Range("S5").Value = 0.000009848620384544656
$ws.Range("T5").Value = 0.000009848620384544656
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Il1rn"
$ws.Range("C6").Value = "Il1r1"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.4115716666666667
$ws.Range("H6").Value = 1.234715
$ws.Range("I6").Value = 0.002110162210096788
$ws.Range("J6").Value = 0.002110162210096788
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 8.314158666666666
$ws.Range("N6").Value = 24.942476
$ws.Range("O6").Value = 0.08379059397414486
$ws.Range("P6").Value = 0.08379059397414489
$ws.Range("Q6").Value = 3.421872139371111
$ws.Range("R6").Value = 30.79684925434
$ws.Range("S6").Value = 0.0001768117449658041
$ws.Range("T6").Value = 0.0001768117449658042
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Il1rn"
$ws.Range("C7").Value = "Il1r1"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 90.57905466666666
$ws.Range("H7").Value = 271.737164
$ws.Range("I7").Value = 0.4644063565694702
$ws.Range("J7").Value = 0.4644063565694702
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 17.79587766666667
$ws.Range("N7").Value = 53.38763299999999
$ws.Range("O7").Value = 0.1793479316144739
$ws.Range("P7").Value = 0.179347931614474
$ws.Range("Q7").Value = 1611.933776010312
$ws.Range("R7").Value = 14507.40398409281
$ws.Range("S7").Value = 0.08329031947934833
$ws.Range("T7").Value = 0.08329031947934835
$ws.Range("A8").Value = "M1"
$ws.Range("B8").Value = "Il1rn"
$ws.Range("C8").Value = "Il1r1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 90.57905466666666
$ws.Range("H8").Value = 271.737164
$ws.Range("I8").Value = 0.4644063565694702
$ws.Range("J8").Value = 0.4644063565694702
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 72.39518466666667
$ws.Range("N8").Value = 217.185554
$ws.Range("O8").Value = 0.7296030503252249
$ws.Range("P8").Value = 0.7296030503252251
$ws.Range("Q8").Value = 6557.487389525428
$ws.Range("R8").Value = 59017.38650572886
$ws.Range("S8").Value = 0.3388322943435095
$ws.Range("T8").Value = 0.3388322943435096
$ws.Range("A9").Value = "M1"
$ws.Range("B9").Value = "Il1rn"
$ws.Range("C9").Value = "Il1r1"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 90.57905466666666
$ws.Range("H9").Value = 271.737164
$ws.Range("I9").Value = 0.4644063565694702
$ws.Range("J9").Value = 0.4644063565694702
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.257112
$ws.Range("N9").Value = 0.771336
$ws.Range("O9").Value = 0.00259119028895291
$ws.Range("P9").Value = 0.00259119028895291
$ws.Range("Q9").Value = 23.288961903456
$ws.Range("R9").Value = 209.600657131104
$ws.Range("S9").Value = 0.001203365241270813
$ws.Range("T9").Value = 0.001203365241270814
$ws.Range("A10").Value = "M1"
$ws.Range("B10").Value = "Il1rn"
$ws.Range("C10").Value = "Il1r1"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 90.57905466666666
$ws.Range("H10").Value = 271.737164
$ws.Range("I10").Value = 0.4644063565694702
$ws.Range("J10").Value = 0.4644063565694702
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4631083333333333
$ws.Range("N10").Value = 1.389325
$ws.Range("O10").Value = 0.004667233797203165
$ws.Range("P10").Value = 0.004667233797203166
$ws.Range("Q10").Value = 41.94791504158889
$ws.Range("R10").Value = 377.5312353743
$ws.Range("S10").Value = 0.002167493043017015
$ws.Range("T10").Value = 0.002167493043017016
$ws.Range("A11").Value = "M1"
$ws.Range("B11").Value = "Il1rn"
$ws.Range("C11").Value = "Il1r1"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 90.57905466666666
$ws.Range("H11").Value = 271.737164
$ws.Range("I11").Value = 0.4644063565694702
$ws.Range("J11").Value = 0.4644063565694702
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 8.314158666666666
$ws.Range("N11").Value = 24.942476
$ws.Range("O11").Value = 0.08379059397414486
$ws.Range("P11").Value = 0.08379059397414489
$ws.Range("Q11").Value = 753.0886323753404
$ws.Range("R11").Value = 6777.797691378064
$ws.Range("S11").Value = 0.03891288446232442
$ws.Range("T11").Value = 0.03891288446232443
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Il1rn"
$ws.Range("C12").Value = "Il1r1"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 104.0520413333333
$ws.Range("H12").Value = 312.156124
$ws.Range("I12").Value = 0.533483481220433
$ws.Range("J12").Value = 0.533483481220433
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 17.79587766666667
$ws.Range("N12").Value = 53.38763299999999
$ws.Range("O12").Value = 0.1793479316144739
$ws.Range("P12").Value = 0.179347931614474
$ws.Range("Q12").Value = 1851.697398534943
$ws.Range("R12").Value = 16665.27658681449
$ws.Range("S12").Value = 0.0956791589073737
$ws.Range("T12").Value = 0.09567915890737373
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Il1rn"
$ws.Range("C13").Value = "Il1r1"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 104.0520413333333
$ws.Range("H13").Value = 312.156124
$ws.Range("I13").Value = 0.533483481220433
$ws.Range("J13").Value = 0.533483481220433
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 72.39518466666667
$ws.Range("N13").Value = 217.185554
$ws.Range("O13").Value = 0.7296030503252249
$ws.Range("P13").Value = 0.7296030503252251
$ws.Range("Q13").Value = 7532.866747270298
$ws.Range("R13").Value = 67795.80072543269
$ws.Range("S13").Value = 0.3892311751965477
$ws.Range("T13").Value = 0.3892311751965479
$ws.Range("A14").Value = "M2"
$ws.Range("B14").Value = "Il1rn"
$ws.Range("C14").Value = "Il1r1"
$ws.Range("D14").Value = "M1"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 104.0520413333333
$ws.Range("H14").Value = 312.156124
$ws.Range("I14").Value = 0.533483481220433
$ws.Range("J14").Value = 0.533483481220433
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.257112
$ws.Range("N14").Value = 0.771336
$ws.Range("O14").Value = 0.00259119028895291
$ws.Range("P14").Value = 0.00259119028895291
$ws.Range("Q14").Value = 26.753028451296
$ws.Range("R14").Value = 240.777256061664
$ws.Range("S14").Value = 0.001382357215855178
$ws.Range("T14").Value = 0.001382357215855178
$ws.Range("A15").Value = "M2"
$ws.Range("B15").Value = "Il1rn"
$ws.Range("C15").Value = "Il1r1"
$ws.Range("D15").Value = "M2"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 104.0520413333333
$ws.Range("H15").Value = 312.156124
$ws.Range("I15").Value = 0.533483481220433
$ws.Range("J15").Value = 0.533483481220433
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.4631083333333333
$ws.Range("N15").Value = 1.389325
$ws.Range("O15").Value = 0.004667233797203165
$ws.Range("P15").Value = 0.004667233797203166
$ws.Range("Q15").Value = 48.1873674418111
$ws.Range("R15").Value = 433.6863069762999
$ws.Range("S15").Value = 0.002489892133801605
$ws.Range("T15").Value = 0.002489892133801605
$ws.Range("A16").Value = "M2"
$ws.Range("B16").Value = "Il1rn"
$ws.Range("C16").Value = "Il1r1"
$ws.Range("D16").Value = "sCs"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 104.0520413333333
$ws.Range("H16").Value = 312.156124
$ws.Range("I16").Value = 0.533483481220433
$ws.Range("J16").Value = 0.533483481220433
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 8.314158666666666
$ws.Range("N16").Value = 24.942476
$ws.Range("O16").Value = 0.08379059397414486
$ws.Range("P16").Value = 0.08379059397414489
$ws.Range("Q16").Value = 865.1051812358913
$ws.Range("R16").Value = 7785.946631123024
$ws.Range("S16").Value = 0.04470089776685463
$ws.Range("T16").Value = 0.04470089776685465
